$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 143
$ws.Range("C3").Value = 143.4
$ws.Range("D3").Value = 144.8
$ws.Range("E3").Value = 145.2
$ws.Range("F3").Value = 145.9
$ws.Range("G3").Value = 148.629146127722
$ws.Range("H3").Value = 147.843274508866
$ws.Range("I3").Value = 159.196596921229
$ws.Range("J3").Value = 170.768954252019
$ws.Range("K3").Value = 182.563822209747
$ws.Range("L3").Value = 193.495476679954
$ws.Range("M3").Value = 213.891482041766
$ws.Range("N3").Value = 233.327403374037
$ws.Range("O3").Value = 252.904562984213
$ws.Range("P3").Value = 272.471324741396
$ws.Range("Q3").Value = 280.386022004872
$ws.Range("R3").Value = 288.71641952255
$ws.Range("S3").Value = 297.424726362232
$ws.Range("T3").Value = 305.603960151109
$ws.Range("U3").Value = 306.706505161232
$ws.Range("V3").Value = 307.771259239154
$ws.Range("W3").Value = 308.798222384877
$ws.Range("X3").Value = 309.825185530599
$ws.Range("Y3").Value = 310.896273431688
$ws.Range("Z3").Value = 311.967361332776
$ws.Range("AA3").Value = 313.076240166065
$ws.Range("AB3").Value = 314.185118999354
$ws.Range("AC3").Value = 313.778526206986
$ws.Range("AD3").Value = 313.371933414617
$ws.Range("AE3").Value = 313.003131554449
$ws.Range("AF3").Value = 312.672120626481
$ws.Range("AG3").Value = 312.242540688541
$ws.Range("AH3").Value = 311.812960750601
$ws.Range("AI3").Value = 311.421171744862
$ws.Range("AJ3").Value = 311.029382739122
$ws.Range("AK3").Value = 310.267230678444
$ws.Range("AL3").Value = 309.505078617766
$ws.Range("AM3").Value = 308.742926557087
$ws.Range("AN3").Value = 308.018565428609
$ws.Range("AO3").Value = 306.671974509099
$ws.Range("AP3").Value = 305.363174521791
$ws.Range("AQ3").Value = 304.129956398882
$ws.Range("AR3").Value = 302.972320140374
$ws.Range("AS3").Value = 301.367885577897
$ws.Range("AT3").Value = 299.83903287982
$ws.Range("AU3").Value = 298.423552978344
$ws.Range("AV3").Value = 297.045864009069